# Weekly update of the "Pimiento" (Macroferia Regional de Talca) sheet:
# a new daily record is inserted as row 312, pushing every existing
# record from row 312 downward by one row (old row 355 becomes row 356).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 312, shifting rows 312:355 down to 313:356.
$ws.Rows.Item(312).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A312").Value = 5
$ws.Range("B312").Value = "Macroferia Regional de Talca"
$ws.Range("C312").Value = "Maule"
$ws.Range("D312").Value = 44491
$ws.Range("E312").Value = 7
$ws.Range("F312").Value = 100112002
$ws.Range("G312").Value = "Pimiento"
$ws.Range("H312").Value = "Zafiro rojo"
$ws.Range("I312").Value = "Primera"
$ws.Range("J312").Value = 200
$ws.Range("K312").Value = 43000
$ws.Range("L312").Value = 43000
$ws.Range("M312").Value = 43000
$ws.Range("N312").Value = "`$/caja 15 kilos"
$ws.Range("O312").Value = "Región de Arica y Parinacota"
$ws.Range("P312").Value = 2867
$ws.Range("Q312").Value = 15
$ws.Range("R312").Value = "Hortaliza"
